$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.03271624445915222
$ws.Range("C2").Value = 0.014678911305963993
$ws.Range("D2").Value = 0.010618066415190697
$ws.Range("E2").Value = 0.011015959084033966
$ws.Range("F2").Value = 0.00007989072037162259
$ws.Range("J2").Value = 0.12764394283294678
$ws.Range("K2").Value = 1.4613769054412842
